$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cell B16 already carries the "answer" cell style (s=3) used throughout
# column B of the grid; copy its formatting onto every B cell that is being
# filled in below so the new cells match the rest of the sheet.
$fmtSource = $ws.Range("B16")

$rowsToFormat = @(147,148,149,151,152,153,154,155,156,157,158,159,160,161,162,163,164,166)
foreach ($r in $rowsToFormat) {
    $fmtSource.Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
}

# Simple "YES" answers
$ws.Cells.Item(147, 2).Value = "YES"
$ws.Cells.Item(148, 2).Value = "YES"
$ws.Cells.Item(149, 2).Value = "YES"
$ws.Cells.Item(151, 2).Value = "YES"
$ws.Cells.Item(152, 2).Value = "YES"
$ws.Cells.Item(153, 2).Value = "YES"
$ws.Cells.Item(154, 2).Value = "YES"
$ws.Cells.Item(155, 2).Value = "YES"
$ws.Cells.Item(156, 2).Value = "YES"
$ws.Cells.Item(157, 2).Value = "YES"
$ws.Cells.Item(158, 2).Value = "YES"
$ws.Cells.Item(159, 2).Value = "YES"
$ws.Cells.Item(160, 2).Value = "YES"
$ws.Cells.Item(162, 2).Value = "YES"
$ws.Cells.Item(164, 2).Value = "YES"

# Commentary / special answers
$ws.Cells.Item(161, 2).Value = "This isn't a yes or no question" + [char]10
$ws.Cells.Item(163, 2).Value = "YES" + [char]10
$ws.Cells.Item(166, 2).Value = "More thatn minimum giant monster, multi combat,  limited players playing simultaneously (birdsong), more than requried treasures/group characters."

# Embedded newlines in B161/B163 otherwise trigger row autofit; restore the
# original explicit row heights so only the cell contents/formatting change.
$ws.Rows.Item(161).RowHeight = 15.75
$ws.Rows.Item(163).RowHeight = 15.75
